$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Lower Right Cell" references (column D, rows 5-11) from row 36 to row 39
# to account for the extra scenarios (40,41,42,44) added to the init file.
$ws.Range("D5").Value  = "A39"
$ws.Range("D6").Value  = "B39"
$ws.Range("D7").Value  = "C39"
$ws.Range("D8").Value  = "G39"
$ws.Range("D9").Value  = "H39"
$ws.Range("D10").Value = "I39"
$ws.Range("D11").Value = "J39"
